# ----------------------------------------------------------------------
# Applies the "Add files via upload" commit to iTTL.docx:
#   1. "Copying" -> "Copying " + "with optional sum" (two runs; the
#      paragraph mark rPr picks up the sz/szCs 20 run formatting).
#   2. Five new paragraphs are inserted right after the
#      "// for each row of M copy V(row_number) ..." example and before
#      "Multiplication": a new numbered list item ("Copying with type
#      conversion without sum"), an "Example:" caption and three
#      monospace code/comment lines demonstrating MF(I,J)=MD(J,I).
# ----------------------------------------------------------------------

$d = $word.ActiveDocument

function Insert-ParagraphXml {
    param($AfterParagraph, $InnerXml)
    $AfterParagraph.Range.InsertParagraphAfter()
    $newPara = $AfterParagraph.Next()
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $InnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newPara.Range.InsertXML($pkg)
    return $newPara
}

# ------------------------------------------------------------------
# 1. "Copying" -> "Copying " + "with optional sum"
# ------------------------------------------------------------------

$copyingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext.Trim() -eq "Copying") {
        $copyingPara = $d.Paragraphs.Item($i)
        break
    }
}

$copyingXml = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>' +
    '<w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Copying </w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>with optional sum</w:t></w:r></w:p>'

$pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $copyingXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$copyingPara.Range.InsertXML($pkg)

# ------------------------------------------------------------------
# 2. Insert the five new paragraphs after the "for each row of M copy
#    V(row_number) ..." comment line and before "Multiplication".
# ------------------------------------------------------------------

$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*for each row of M copy V(row_number)*") {
        $anchor = $d.Paragraphs.Item($i)
        break
    }
}

# -- 2a. "Copying with type conversion without sum" (new numbered item)
$xml1 = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr/></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Copying with type conversion without sum</w:t></w:r></w:p>'
$anchor = Insert-ParagraphXml $anchor $xml1

# -- 2b. "Example:"
$xml2 = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
    '<w:ind w:left="720" w:hanging="0"/><w:rPr/></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Example:</w:t></w:r></w:p>'
$anchor = Insert-ParagraphXml $anchor $xml2

# -- 2c. iTTL::MATRIX<float> MF({2,2});  // uninitialized matrix 2x2 of float values
$xml3 = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr>' +
    '<w:rFonts w:ascii="Monospace" w:hAnsi="Monospace" w:eastAsia="Noto Sans CJK SC Regular" w:cs="Lohit Devanagari"/>' +
    '<w:color w:val="000000"/><w:kern w:val="2"/><w:sz w:val="16"/><w:szCs w:val="16"/>' +
    '<w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Lohit Devanagari" w:ascii="Monospace" w:hAnsi="Monospace"/>' +
    '<w:color w:val="000000"/><w:kern w:val="2"/><w:sz w:val="16"/><w:szCs w:val="16"/>' +
    '<w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr>' +
    '<w:t xml:space="preserve">iTTL::MATRIX&lt;float&gt; MF({2,2}); </w:t><w:tab/><w:tab/><w:tab/></w:r>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Lohit Devanagari" w:ascii="Monospace" w:hAnsi="Monospace"/>' +
    '<w:color w:val="3F7F5F"/><w:kern w:val="2"/><w:sz w:val="16"/><w:szCs w:val="16"/>' +
    '<w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr>' +
    '<w:t>// uninitialized matrix 2x2 of float values</w:t></w:r></w:p>'
$anchor = Insert-ParagraphXml $anchor $xml3

# -- 2d. iTTL::MATRIX<double> MD({{1.0,2.0},{3.0,4.0}});  // initialized matrix 2x2 of double values
$xml4 = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr>' +
    '<w:rFonts w:ascii="Monospace" w:hAnsi="Monospace" w:eastAsia="Noto Sans CJK SC Regular" w:cs="Lohit Devanagari"/>' +
    '<w:color w:val="000000"/><w:kern w:val="2"/><w:sz w:val="16"/><w:szCs w:val="16"/>' +
    '<w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Lohit Devanagari" w:ascii="Monospace" w:hAnsi="Monospace"/>' +
    '<w:color w:val="000000"/><w:kern w:val="2"/><w:sz w:val="16"/><w:szCs w:val="16"/>' +
    '<w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr>' +
    '<w:t xml:space="preserve">iTTL::MATRIX&lt;double&gt; MD({{1.0,2.0},{3.0,4.0}}); </w:t><w:tab/></w:r>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Lohit Devanagari" w:ascii="Monospace" w:hAnsi="Monospace"/>' +
    '<w:color w:val="3F7F5F"/><w:kern w:val="2"/><w:sz w:val="16"/><w:szCs w:val="16"/>' +
    '<w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr>' +
    '<w:t>// initialized matrix 2x2 of double values</w:t></w:r></w:p>'
$anchor = Insert-ParagraphXml $anchor $xml4

# -- 2e. MF(I,J)=MD(J,I);  // copy with transposition
$xml5 = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr>' +
    '<w:rFonts w:ascii="Monospace" w:hAnsi="Monospace" w:eastAsia="Noto Sans CJK SC Regular" w:cs="Lohit Devanagari"/>' +
    '<w:color w:val="000000"/><w:kern w:val="2"/><w:sz w:val="16"/><w:szCs w:val="16"/>' +
    '<w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Lohit Devanagari" w:ascii="Monospace" w:hAnsi="Monospace"/>' +
    '<w:color w:val="000000"/><w:kern w:val="2"/><w:sz w:val="16"/><w:szCs w:val="16"/>' +
    '<w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr>' +
    '<w:t xml:space="preserve">MF(I,J)=MD(J,I); </w:t><w:tab/><w:tab/><w:tab/><w:tab/><w:tab/></w:r>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Lohit Devanagari" w:ascii="Monospace" w:hAnsi="Monospace"/>' +
    '<w:color w:val="3F7F5F"/><w:kern w:val="2"/><w:sz w:val="16"/><w:szCs w:val="16"/>' +
    '<w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr>' +
    '<w:t>// copy with transposition</w:t></w:r></w:p>'
$anchor = Insert-ParagraphXml $anchor $xml5

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
